$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.746.85"
$ws.Range("E2").Value = "  +1.40%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.567.47"
$ws.Range("E3").Value = "  -0.37%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.38%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.13"
$ws.Range("E5").Value = "  -0.83%  "

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.488"
$ws.Range("E6").Value = "  -0.24%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.28%  "

# Row 8 - Solana
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "25.32"
$ws.Range("E8").Value = "  +6.53%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.246"
$ws.Range("E9").Value = "  +0.01%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0586"
$ws.Range("E10").Value = "  -0.22%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0894"
$ws.Range("E11").Value = "  -0.35%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.787.47"
$ws.Range("E12").Value = "  -0.65%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.565.42"
$ws.Range("E13").Value = "  -0.67%  "

# Row 14 - WrappedBTC
$ws.Range("D14").Value = "28.728.13"
$ws.Range("E14").Value = "  +1.28%  "

# Row 15 - Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.517"
$ws.Range("E15").Value = "  +0.11%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  -1.24%  "

# Row 17 - Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.57"
$ws.Range("E17").Value = "  +0.02%  "

# Row 18 - BitcoinCash
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.75"
$ws.Range("E18").Value = "  +1.42%  "

# Row 19 - Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.39"
$ws.Range("E19").Value = "  -0.51%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0676"
$ws.Range("E20").Value = "  -1.07%  "

# Row 21 - Dai
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.997"
$ws.Range("E21").Value = "  -0.39%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.93"
$ws.Range("E22").Value = "  -0.99%  "

# Row 23 - Avalanche
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.06"
$ws.Range("E23").Value = "  +0.32%  "

# Row 24 - Toncoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.10"
$ws.Range("E24").Value = "  +2.98%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.80"
$ws.Range("E25").Value = "  -0.33%  "

# Row 26 - EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.84"
$ws.Range("E26").Value = "  -0.59%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.10%  "

# Row 28 - BinanceUSD
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.37%  "

# Row 29 - Cosmos
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.25"
$ws.Range("E29").Value = "  -1.83%  "

# Row 30 - Hedera
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0462"
$ws.Range("E30").Value = "  -3.78%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.56%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.73%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.389.10"
$ws.Range("E33").Value = "  +0.23%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.99"
$ws.Range("E34").Value = "  -3.50%  "

# Row 35 - TrustWalletToken
$ws.Range("E35").Value = "  -2.79%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  -1.75%  "

# Row 37 - now MXToken (was HuobiToken)
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.69"
$ws.Range("E37").Value = "  +2.11%  "

# Row 38 - now HuobiToken (was MXToken)
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.30"
$ws.Range("E38").Value = "  -2.71%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -0.96%  "

# Row 40 - RenderToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.96"
$ws.Range("E40").Value = "  +3.06%  "

# Row 41 - ImmutableX
$ws.Range("E41").Value = "  +0.38%  "

# Row 42 - PaxDollar
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  -0.28%  "

# Row 43 - ARBITRUM
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.777"
$ws.Range("E43").Value = "  -1.19%  "

# Row 44 - Kaspa
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0460"
$ws.Range("E44").Value = "  -1.18%  "

# Row 45 - Aave
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.27"
$ws.Range("E45").Value = "  +3.08%  "

# Row 46 - FraxShare
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.29"
$ws.Range("E46").Value = "  -1.73%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.702.65"
$ws.Range("E47").Value = "  -0.51%  "

# Row 48 - WEMIXToken
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.868"
$ws.Range("E48").Value = "  -5.57%  "

# Row 49 - now BitcoinSV (was Quant)
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "43.92"
$ws.Range("E49").Value = "  +6.85%  "

# Row 50 - now Quant (was BitcoinSV)
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.41"
$ws.Range("E50").Value = "  +0.13%  "

# Row 51 - now Cronos (was BabyDogeCoin)
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0509"
$ws.Range("E51").Value = "  -1.04%  "
